$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D-column cells whose new values would otherwise
# be auto-converted to numbers by Excel (losing exact text representation,
# e.g. trailing zeros such as "0.550").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.697.50"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.607.11"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.79%  "
$ws.Range("D5").Value = "212.94"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "0.517"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").Value = "28.11"
$ws.Range("E8").Value = "  +5.18%  "
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").Value = "0.0603"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.837.49"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "1.604.85"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "0.548"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").Value = "29.713.28"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "3.76"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "64.17"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "241.56"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "7.85"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").Value = "155.19"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").Value = "0.0481"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").Value = "1.428.32"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.550"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").Value = "56.87"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("E42").Value = "  +5.95%  "
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "66.36"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").Value = "0.982"
$ws.Range("E47").Value = "  +17.57%  "
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").Value = "1.746.44"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "86.66"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("E51").Value = "  +5.32%  "
